$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 147 - this shifts the existing rows 147 and 148
# down to 148 and 149 respectively (with all their values/formatting intact),
# and grows the used range to A1:R149.
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new weekly price record.
$ws.Cells.Item(147, 1).Value = 10
$ws.Cells.Item(147, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(147, 3).Value = "La Araucanía"
$ws.Cells.Item(147, 4).Value = 44448
$ws.Cells.Item(147, 5).Value = 9
$ws.Cells.Item(147, 6).Value = 100112043
$ws.Cells.Item(147, 7).Value = "Pepino dulce"
$ws.Cells.Item(147, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 150
$ws.Cells.Item(147, 11).Value = 22000
$ws.Cells.Item(147, 12).Value = 22000
$ws.Cells.Item(147, 13).Value = 22000
$ws.Cells.Item(147, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(147, 15).Value = "Provincia de Copiapó"
$ws.Cells.Item(147, 16).Value = 1222
$ws.Cells.Item(147, 17).Value = 18
$ws.Cells.Item(147, 18).Value = "Hortaliza"
